$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '62.071.31'
$ws.Range('E2').Value = '  +0.01%  '
$ws.Range('D3').Value = '3.428.88'
$ws.Range('E3').Value = '  +0.15%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '409.61'
$ws.Range('D5').NumberFormat = 'General'
$ws.Range('E5').Value = '  +0.70%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '130.22'
$ws.Range('D6').NumberFormat = 'General'
$ws.Range('E6').Value = '  -0.97%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.635'
$ws.Range('D7').NumberFormat = 'General'
$ws.Range('E7').Value = '  +6.97%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '1.00'
$ws.Range('D8').NumberFormat = 'General'
$ws.Range('E8').Value = '  -0.02%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.741'
$ws.Range('D9').NumberFormat = 'General'
$ws.Range('E9').Value = '  +7.27%  '
$ws.Range('E10').Value = '  +5.71%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '42.89'
$ws.Range('D11').NumberFormat = 'General'
$ws.Range('E11').Value = '  +2.46%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0000225'
$ws.Range('D12').NumberFormat = 'General'
$ws.Range('E12').Value = '  +51.70%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '9.25'
$ws.Range('D13').NumberFormat = 'General'
$ws.Range('E13').Value = '  +10.60%  '
$ws.Range('E14').Value = '  -0.23%  '
$ws.Range('D15').Value = '3.973.93'
$ws.Range('E15').Value = '  -0.28%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '21.34'
$ws.Range('D16').NumberFormat = 'General'
$ws.Range('E16').Value = '  +7.79%  '
$ws.Range('D17').Value = '3.422.13'
$ws.Range('E17').Value = '  -1.24%  '
$ws.Range('E18').Value = '  +7.81%  '
$ws.Range('E19').Value = '  +8.16%  '
$ws.Range('D20').Value = '62.035.07'
$ws.Range('E20').Value = '  +0.03%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '454.52'
$ws.Range('D21').NumberFormat = 'General'
$ws.Range('E21').Value = '  +45.77%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '91.49'
$ws.Range('D22').NumberFormat = 'General'
$ws.Range('E22').Value = '  +9.24%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '3.22'
$ws.Range('D23').NumberFormat = 'General'
$ws.Range('E23').Value = '  +1.63%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '13.14'
$ws.Range('D24').NumberFormat = 'General'
$ws.Range('E24').Value = '  +2.91%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '3.28'
$ws.Range('D25').NumberFormat = 'General'
$ws.Range('E25').Value = '  +3.43%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '33.17'
$ws.Range('D26').NumberFormat = 'General'
$ws.Range('E26').Value = '  +11.93%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.12'
$ws.Range('D27').NumberFormat = 'General'
$ws.Range('E27').Value = '  +11.91%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '4.78'
$ws.Range('D28').NumberFormat = 'General'
$ws.Range('E28').Value = '  +1.63%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.66'
$ws.Range('D29').NumberFormat = 'General'
$ws.Range('E29').Value = '  -1.55%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '12.11'
$ws.Range('D30').NumberFormat = 'General'
$ws.Range('E30').Value = '  +6.65%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '2.68'
$ws.Range('D31').NumberFormat = 'General'
$ws.Range('E31').Value = '  -2.27%  '
$ws.Range('E32').Value = '  -0.50%  '
$ws.Range('E33').Value = '  -0.16%  '
$ws.Range('E34').Value = '  -1.29%  '
$ws.Range('E35').Value = '  -0.10%  '
$ws.Range('E36').Value = '  +3.68%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '54.32'
$ws.Range('D37').NumberFormat = 'General'
$ws.Range('E37').Value = '  +4.89%  '
$ws.Range('E38').Value = '  -0.06%  '
$ws.Range('E39').Value = '  +2.14%  '
$ws.Range('E40').Value = '  +7.78%  '
$ws.Range('E41').Value = '  -1.01%  '
$ws.Range('E42').Value = '  +0.39%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '142.55'
$ws.Range('D43').NumberFormat = 'General'
$ws.Range('E43').Value = '  -0.50%  '
$ws.Range('E44').Value = '  +9.80%  '
$ws.Range('E45').Value = '  +1.54%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.54'
$ws.Range('D46').NumberFormat = 'General'
$ws.Range('E46').Value = '  +14.72%  '
$ws.Range('E47').Value = '  -0.93%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '22.50'
$ws.Range('D48').NumberFormat = 'General'
$ws.Range('E48').Value = '  +5.84%  '
$ws.Range('B49').Value = 'Cronos'
$ws.Range('C49').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.139'
$ws.Range('D49').NumberFormat = 'General'
$ws.Range('E49').Value = '  +17.17%  '
$ws.Range('B50').Value = 'RocketPoolETH'
$ws.Range('C50').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D50').Value = '3.778.08'
$ws.Range('E50').Value = '  -0.11%  '
$ws.Range('B51').Value = 'ThetaToken'
$ws.Range('C51').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.12'
$ws.Range('D51').NumberFormat = 'General'
$ws.Range('E51').Value = '  +8.69%  '
